$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values: column E (times) ---
$ws.Range("E4").Value = 0.015
$ws.Range("E5").Value = 0.096
$ws.Range("E6").Value = 0.654
$ws.Range("E7").Value = 2.762
$ws.Range("E8").Value = 11.53
$ws.Range("E9").Value = 53.461
$ws.Range("E10").Value = 78.015

# --- D10: replace formula with literal value (shrinks the shared-formula range D6:D11 -> D6:D10) ---
$ws.Range("D10").Value = 10000

# --- Number format for time column ---
$ws.Range("E4:E10").NumberFormat = "0.000"

# --- Column E width ---
$ws.Range("E:E").ColumnWidth = 10.9

# --- Borders: thin box border around D3:E10 cells individually ---
$dataRange = $ws.Range("D3:E10")
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# --- Fills: header row D3:E3 and data column D4:D10 ---
$headerRange = $ws.Range("D3:E3")
$headerRange.Interior.ThemeColor = 10
$headerRange.Interior.TintAndShade = 0.59999389629810485

$col1Range = $ws.Range("D4:D10")
$col1Range.Interior.ThemeColor = 9
$col1Range.Interior.TintAndShade = 0.59999389629810485

# --- Selection ---
$ws.Range("D3:E10").Select()
